# Search.xlsx - "Test Cases" sheet: Runmode (column D) bulk update.
#
# Row 2 (the first test case) keeps Runmode = "Y"; every other test case row
# (rows 3 through 127) is switched from "Y" to "N" so only TCID row 2 still
# runs. This mirrors an Excel user typing "N" in D3 and filling it down to
# D127, then leaving that range selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("D3:D127").Value = "N"

# Leave the same range selected, matching the post-edit selection state.
$ws.Range("D3:D127").Select()
